$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 12 - A/B/C columns (Year, Month, Day)
$ws.Cells.Item(12, 1).Value = 2026
$ws.Cells.Item(12, 2).Value = 2
$ws.Cells.Item(12, 3).Value = 15

$vals7 = @(4980,4731,4615,4583,4605,4745,5257,5838,6343,6561,6517,6447,6379,6338,6275,6267,6289,6409,6610,6650,6350,6087,5853,5508)
for ($i = 0; $i -lt $vals7.Length; $i++) {
    $ws.Cells.Item(7, 4 + $i).Value = $vals7[$i]
}

$vals8 = @(5121,4872,4756,4745,4724,4876,5354,5896,6368,6571,6527,6457,6388,6346,6283,6274,6284,6332,6431,6398,6181,5883,5655,5319)
for ($i = 0; $i -lt $vals8.Length; $i++) {
    $ws.Cells.Item(8, 4 + $i).Value = $vals8[$i]
}

$vals9 = @(4947,4697,4581,4570,4548,4700,5178,5720,6192,6395,6359,6302,6246,6211,6161,6154,6160,6193,6260,6227,6008,5707,5477,5138)
for ($i = 0; $i -lt $vals9.Length; $i++) {
    $ws.Cells.Item(9, 4 + $i).Value = $vals9[$i]
}

$vals10 = @(4798,4548,4431,4420,4398,4550,5027,5569,6041,6244,6215,6169,6124,6097,6056,6050,6050,6044,6034,6001,5787,5492,5267,4935)
for ($i = 0; $i -lt $vals10.Length; $i++) {
    $ws.Cells.Item(10, 4 + $i).Value = $vals10[$i]
}

$vals11 = @(4597,4371,4261,4143,4172,4305,4493,4725,5225,5536,5737,5721,5715,5696,5682,5675,5672,5656,5641,5496,5256,4967,4871,4636)
for ($i = 0; $i -lt $vals11.Length; $i++) {
    $ws.Cells.Item(11, 4 + $i).Value = $vals11[$i]
}

$vals12 = @(4512,4256,4109,4067,4106,4192,4402,4652,5074,5380,5700,5705,5706,5707,5711,5720,5755,5800,5825,5829,5536,5268,5066,4699)
for ($i = 0; $i -lt $vals12.Length; $i++) {
    $ws.Cells.Item(12, 4 + $i).Value = $vals12[$i]
}

